$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8966636666666666
$ws.Range("H2").Value = 2.689991
$ws.Range("I2").Value = 0.334725143386341
$ws.Range("J2").Value = 0.3647360854412732
$ws.Range("O2").Value = 0.05741109988842188
$ws.Range("P2").Value = 0.05868788750479208
$ws.Range("Q2").Value = 0.2362813372427778
$ws.Range("R2").Value = 2.126532035185
$ws.Range("S2").Value = 0.01921693864211956
$ws.Range("T2").Value = 0.02140559035131567
$ws.Range("G3").Value = 0.8966636666666666
$ws.Range("H3").Value = 2.689991
$ws.Range("I3").Value = 0.334725143386341
$ws.Range("J3").Value = 0.3647360854412732
$ws.Range("M3").Value = 2.338082
$ws.Range("N3").Value = 7.014246
$ws.Range("O3").Value = 0.5093962667661314
$ws.Range("P3").Value = 0.5207249270164355
$ws.Range("Q3").Value = 2.096473179087333
$ws.Range("R3").Value = 18.868258611786
$ws.Range("S3").Value = 0.1705077384337601
$ws.Range("T3").Value = 0.1899271714716674
$ws.Range("G4").Value = 0.8966636666666666
$ws.Range("H4").Value = 2.689991
$ws.Range("I4").Value = 0.334725143386341
$ws.Range("J4").Value = 0.3647360854412732
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2995679999999999
$ws.Range("N4").Value = 0.5991359999999999
$ws.Range("O4").Value = 0.06526666765434078
$ws.Range("P4").Value = 0.04447877218348473
$ws.Range("Q4").Value = 0.268611741296
$ws.Range("R4").Value = 1.611670447776
$ws.Range("S4").Value = 0.02184639468894788
$ws.Range("T4").Value = 0.01622301325143841
$ws.Range("G5").Value = 0.8966636666666666
$ws.Range("H5").Value = 2.689991
$ws.Range("I5").Value = 0.334725143386341
$ws.Range("J5").Value = 0.3647360854412732
$ws.Range("M5").Value = 1.688746333333333
$ws.Range("N5").Value = 5.066238999999999
$ws.Range("O5").Value = 0.3679259656911061
$ws.Range("P5").Value = 0.3761084132952877
$ws.Range("Q5").Value = 1.514237479316555
$ws.Range("R5").Value = 13.628137313849
$ws.Range("S5").Value = 0.1231540716215135
$ws.Range("T5").Value = 0.1371803103668518
$ws.Range("G6").Value = 0.6827986666666667
$ws.Range("I6").Value = 0.2548891965854188
$ws.Range("J6").Value = 0.2777421703171357
$ws.Range("O6").Value = 0.05741109988842188
$ws.Range("P6").Value = 0.05868788750479208
$ws.Range("S6").Value = 0.01463346912564508
$ws.Range("T6").Value = 0.01630010124690886
$ws.Range("G7").Value = 0.6827986666666667
$ws.Range("I7").Value = 0.2548891965854188
$ws.Range("J7").Value = 0.2777421703171357
$ws.Range("M7").Value = 2.338082
$ws.Range("N7").Value = 7.014246
$ws.Range("O7").Value = 0.5093962667661314
$ws.Range("P7").Value = 0.5207249270164355
$ws.Range("Q7").Value = 1.596439272157333
$ws.Range("R7").Value = 14.367953449416
$ws.Range("S7").Value = 0.1298396051796309
$ws.Range("T7").Value = 0.1446272713677769
$ws.Range("G8").Value = 0.6827986666666667
$ws.Range("I8").Value = 0.2548891965854188
$ws.Range("J8").Value = 0.2777421703171357
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2995679999999999
$ws.Range("N8").Value = 0.5991359999999999
$ws.Range("O8").Value = 0.06526666765434078
$ws.Range("P8").Value = 0.04447877218348473
$ws.Range("Q8").Value = 0.204544630976
$ws.Range("R8").Value = 1.227267785856
$ws.Range("S8").Value = 0.01663576848222246
$ws.Range("T8").Value = 0.01235363071928249
$ws.Range("G9").Value = 0.6827986666666667
$ws.Range("I9").Value = 0.2548891965854188
$ws.Range("J9").Value = 0.2777421703171357
$ws.Range("M9").Value = 1.688746333333333
$ws.Range("N9").Value = 5.066238999999999
$ws.Range("O9").Value = 0.3679259656911061
$ws.Range("P9").Value = 0.3761084132952877
$ws.Range("Q9").Value = 1.153073744738222
$ws.Range("R9").Value = 10.377663702644
$ws.Range("S9").Value = 0.0937803537979204
$ws.Range("T9").Value = 0.1044611669831675
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.221369
$ws.Range("H10").Value = 0.664107
$ws.Range("I10").Value = 0.08263719499391366
$ws.Range("J10").Value = 0.09004631892602898
$ws.Range("O10").Value = 0.05741109988842188
$ws.Range("P10").Value = 0.05868788750479208
$ws.Range("Q10").Value = 0.05833331413833333
$ws.Range("R10").Value = 0.5249998272449999
$ws.Range("S10").Value = 0.004744292256294573
$ws.Range("T10").Value = 0.005284628235351418
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.221369
$ws.Range("H11").Value = 0.664107
$ws.Range("I11").Value = 0.08263719499391366
$ws.Range("J11").Value = 0.09004631892602898
$ws.Range("M11").Value = 2.338082
$ws.Range("N11").Value = 7.014246
$ws.Range("O11").Value = 0.5093962667661314
$ws.Range("P11").Value = 0.5207249270164355
$ws.Range("Q11").Value = 0.517578874258
$ws.Range("R11").Value = 4.658209868322
$ws.Range("S11").Value = 0.04209507862592446
$ws.Range("T11").Value = 0.04688936285085511
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.221369
$ws.Range("H12").Value = 0.664107
$ws.Range("I12").Value = 0.08263719499391366
$ws.Range("J12").Value = 0.09004631892602898
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2995679999999999
$ws.Range("N12").Value = 0.5991359999999999
$ws.Range("O12").Value = 0.06526666765434078
$ws.Range("P12").Value = 0.04447877218348473
$ws.Range("Q12").Value = 0.066315068592
$ws.Range("R12").Value = 0.3978904115519999
$ws.Range("S12").Value = 0.005393454341554716
$ws.Range("T12").Value = 0.004005149705472252
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.221369
$ws.Range("H13").Value = 0.664107
$ws.Range("I13").Value = 0.08263719499391366
$ws.Range("J13").Value = 0.09004631892602898
$ws.Range("M13").Value = 1.688746333333333
$ws.Range("N13").Value = 5.066238999999999
$ws.Range("O13").Value = 0.3679259656911061
$ws.Range("P13").Value = 0.3761084132952877
$ws.Range("Q13").Value = 0.3738360870636667
$ws.Range("R13").Value = 3.364524783573
$ws.Range("S13").Value = 0.03040436977013992
$ws.Range("T13").Value = 0.0338671781343502
$ws.Range("G14").Value = 0.6612465000000001
$ws.Range("H14").Value = 1.322493
$ws.Range("I14").Value = 0.2468437584284291
$ws.Range("J14").Value = 0.1793169270244717
$ws.Range("O14").Value = 0.05741109988842188
$ws.Range("P14").Value = 0.05868788750479208
$ws.Range("Q14").Value = 0.1742461672925
$ws.Range("R14").Value = 1.045477003755
$ws.Range("S14").Value = 0.01417157167196802
$ws.Range("T14").Value = 0.01052373164091721
$ws.Range("G15").Value = 0.6612465000000001
$ws.Range("H15").Value = 1.322493
$ws.Range("I15").Value = 0.2468437584284291
$ws.Range("J15").Value = 0.1793169270244717
$ws.Range("M15").Value = 2.338082
$ws.Range("N15").Value = 7.014246
$ws.Range("O15").Value = 0.5093962667661314
$ws.Range("P15").Value = 0.5207249270164355
$ws.Range("Q15").Value = 1.546048539213
$ws.Range("R15").Value = 9.276291235278
$ws.Range("S15").Value = 0.1257412890179626
$ws.Range("T15").Value = 0.09337479373762954
$ws.Range("G16").Value = 0.6612465000000001
$ws.Range("H16").Value = 1.322493
$ws.Range("I16").Value = 0.2468437584284291
$ws.Range("J16").Value = 0.1793169270244717
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2995679999999999
$ws.Range("N16").Value = 0.5991359999999999
$ws.Range("O16").Value = 0.06526666765434078
$ws.Range("P16").Value = 0.04447877218348473
$ws.Range("Q16").Value = 0.198088291512
$ws.Range("R16").Value = 0.792353166048
$ws.Range("S16").Value = 0.01611066954389666
$ws.Range("T16").Value = 0.007975796745764034
$ws.Range("G17").Value = 0.6612465000000001
$ws.Range("H17").Value = 1.322493
$ws.Range("I17").Value = 0.2468437584284291
$ws.Range("J17").Value = 0.1793169270244717
$ws.Range("M17").Value = 1.688746333333333
$ws.Range("N17").Value = 5.066238999999999
$ws.Range("O17").Value = 0.3679259656911061
$ws.Range("P17").Value = 0.3761084132952877
$ws.Range("Q17").Value = 1.1166776023045
$ws.Range("R17").Value = 6.700065613827
$ws.Range("S17").Value = 0.09082022819460189
$ws.Range("T17").Value = 0.06744260490016096
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.216728
$ws.Range("H18").Value = 0.650184
$ws.Range("I18").Value = 0.08090470660589748
$ws.Range("J18").Value = 0.08815849829109049
$ws.Range("O18").Value = 0.05741109988842188
$ws.Range("P18").Value = 0.05868788750479208
$ws.Range("Q18").Value = 0.05711035649333333
$ws.Range("R18").Value = 0.51399320844
$ws.Range("S18").Value = 0.004644828192394646
$ws.Range("T18").Value = 0.005173836030298923
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.216728
$ws.Range("H19").Value = 0.650184
$ws.Range("I19").Value = 0.08090470660589748
$ws.Range("J19").Value = 0.08815849829109049
$ws.Range("M19").Value = 2.338082
$ws.Range("N19").Value = 7.014246
$ws.Range("O19").Value = 0.5093962667661314
$ws.Range("P19").Value = 0.5207249270164355
$ws.Range("Q19").Value = 0.506727835696
$ws.Range("R19").Value = 4.560550521264
$ws.Range("S19").Value = 0.04121255550885334
$ws.Range("T19").Value = 0.04590632758850664
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.216728
$ws.Range("H20").Value = 0.650184
$ws.Range("I20").Value = 0.08090470660589748
$ws.Range("J20").Value = 0.08815849829109049
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.2995679999999999
$ws.Range("N20").Value = 0.5991359999999999
$ws.Range("O20").Value = 0.06526666765434078
$ws.Range("P20").Value = 0.04447877218348473
$ws.Range("Q20").Value = 0.06492477350399999
$ws.Range("R20").Value = 0.3895486410239999
$ws.Range("S20").Value = 0.005280380597719059
$ws.Range("T20").Value = 0.003921181761527542
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.216728
$ws.Range("H21").Value = 0.650184
$ws.Range("I21").Value = 0.08090470660589748
$ws.Range("J21").Value = 0.08815849829109049
$ws.Range("M21").Value = 1.688746333333333
$ws.Range("N21").Value = 5.066238999999999
$ws.Range("O21").Value = 0.3679259656911061
$ws.Range("P21").Value = 0.3761084132952877
$ws.Range("Q21").Value = 0.3659986153306666
$ws.Range("R21").Value = 3.293987537976
$ws.Range("S21").Value = 0.02976694230693044
$ws.Range("T21").Value = 0.03315715291075737
